# Update attendance/view counts (column F) on the "展览" and "全部类型" sheets.
# These counters increased between the previous GitHub Pages data export and
# this one (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 176
$ws1.Range("F3").Value = 662
$ws1.Range("F4").Value = 26
$ws1.Range("F6").Value = 1636
$ws1.Range("F8").Value = 3173
$ws1.Range("F9").Value = 455
$ws1.Range("F10").Value = 750

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 176
$ws4.Range("F3").Value = 662
$ws4.Range("F4").Value = 26
$ws4.Range("F7").Value = 1636
$ws4.Range("F9").Value = 3173
$ws4.Range("F10").Value = 455
$ws4.Range("F11").Value = 750
